$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows at 879, shifting existing rows 879-945 down to 885-951
$ws.Rows("879:884").Insert()

# Row 879: 1a amarillo
$ws.Range("A879").Value = 2
$ws.Range("B879").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C879").Value = "Coquimbo"
$ws.Range("D879").Value = 45147
$ws.Range("E879").Value = 4
$ws.Range("F879").Value = "Fruta"
$ws.Range("G879").Value = 100102
$ws.Range("H879").Value = "Cítricos"
$ws.Range("I879").Value = 100102003
$ws.Range("J879").Value = "Limón"
$ws.Range("K879").Value = "Sin especificar"
$ws.Range("L879").Value = "1a amarillo"
$ws.Range("M879").Value = 420
$ws.Range("N879").Value = 3500
$ws.Range("O879").Value = 3700
$ws.Range("P879").Value = 3600
$ws.Range("Q879").Value = "$/malla 18 kilos"
$ws.Range("R879").Value = "Provincia de Limarí"
$ws.Range("S879").Value = 200
$ws.Range("T879").Value = 18

# Row 880: 1a plateado
$ws.Range("A880").Value = 2
$ws.Range("B880").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C880").Value = "Coquimbo"
$ws.Range("D880").Value = 45147
$ws.Range("E880").Value = 4
$ws.Range("F880").Value = "Fruta"
$ws.Range("G880").Value = 100102
$ws.Range("H880").Value = "Cítricos"
$ws.Range("I880").Value = 100102003
$ws.Range("J880").Value = "Limón"
$ws.Range("K880").Value = "Sin especificar"
$ws.Range("L880").Value = "1a plateado"
$ws.Range("M880").Value = 580
$ws.Range("N880").Value = 3500
$ws.Range("O880").Value = 3700
$ws.Range("P880").Value = 3600
$ws.Range("Q880").Value = "$/malla 18 kilos"
$ws.Range("R880").Value = "Provincia de Limarí"
$ws.Range("S880").Value = 200
$ws.Range("T880").Value = 18

# Row 881: 2a amarillo
$ws.Range("A881").Value = 2
$ws.Range("B881").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C881").Value = "Coquimbo"
$ws.Range("D881").Value = 45147
$ws.Range("E881").Value = 4
$ws.Range("F881").Value = "Fruta"
$ws.Range("G881").Value = 100102
$ws.Range("H881").Value = "Cítricos"
$ws.Range("I881").Value = 100102003
$ws.Range("J881").Value = "Limón"
$ws.Range("K881").Value = "Sin especificar"
$ws.Range("L881").Value = "2a amarillo"
$ws.Range("M881").Value = 360
$ws.Range("N881").Value = 2500
$ws.Range("O881").Value = 2700
$ws.Range("P881").Value = 2600
$ws.Range("Q881").Value = "$/malla 18 kilos"
$ws.Range("R881").Value = "Provincia de Limarí"
$ws.Range("S881").Value = 144
$ws.Range("T881").Value = 18

# Row 882: 2a plateado
$ws.Range("A882").Value = 2
$ws.Range("B882").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C882").Value = "Coquimbo"
$ws.Range("D882").Value = 45147
$ws.Range("E882").Value = 4
$ws.Range("F882").Value = "Fruta"
$ws.Range("G882").Value = 100102
$ws.Range("H882").Value = "Cítricos"
$ws.Range("I882").Value = 100102003
$ws.Range("J882").Value = "Limón"
$ws.Range("K882").Value = "Sin especificar"
$ws.Range("L882").Value = "2a plateado"
$ws.Range("M882").Value = 420
$ws.Range("N882").Value = 2500
$ws.Range("O882").Value = 2700
$ws.Range("P882").Value = 2607
$ws.Range("Q882").Value = "$/malla 18 kilos"
$ws.Range("R882").Value = "Provincia de Limarí"
$ws.Range("S882").Value = 145
$ws.Range("T882").Value = 18

# Row 883: 3a amarillo
$ws.Range("A883").Value = 2
$ws.Range("B883").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C883").Value = "Coquimbo"
$ws.Range("D883").Value = 45147
$ws.Range("E883").Value = 4
$ws.Range("F883").Value = "Fruta"
$ws.Range("G883").Value = 100102
$ws.Range("H883").Value = "Cítricos"
$ws.Range("I883").Value = 100102003
$ws.Range("J883").Value = "Limón"
$ws.Range("K883").Value = "Sin especificar"
$ws.Range("L883").Value = "3a amarillo"
$ws.Range("M883").Value = 300
$ws.Range("N883").Value = 1500
$ws.Range("O883").Value = 1700
$ws.Range("P883").Value = 1600
$ws.Range("Q883").Value = "$/malla 18 kilos"
$ws.Range("R883").Value = "Provincia de Limarí"
$ws.Range("S883").Value = 89
$ws.Range("T883").Value = 18

# Row 884: 3a plateado
$ws.Range("A884").Value = 2
$ws.Range("B884").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C884").Value = "Coquimbo"
$ws.Range("D884").Value = 45147
$ws.Range("E884").Value = 4
$ws.Range("F884").Value = "Fruta"
$ws.Range("G884").Value = 100102
$ws.Range("H884").Value = "Cítricos"
$ws.Range("I884").Value = 100102003
$ws.Range("J884").Value = "Limón"
$ws.Range("K884").Value = "Sin especificar"
$ws.Range("L884").Value = "3a plateado"
$ws.Range("M884").Value = 340
$ws.Range("N884").Value = 1500
$ws.Range("O884").Value = 1700
$ws.Range("P884").Value = 1606
$ws.Range("Q884").Value = "$/malla 18 kilos"
$ws.Range("R884").Value = "Provincia de Limarí"
$ws.Range("S884").Value = 89
$ws.Range("T884").Value = 18
